$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 ---------------------------------------------------------------
$ws.Range("A16").Value = 1261002319
$ws.Range("B16").Value = 7
$ws.Range("D16").Value = "BM"
$ws.Range("F16").Value = "NO"
$ws.Range("Z16").Value = "Desktop"
$ws.Range("N16").Value = "engenharia química bloco 2"
$ws.Range("AA16").Value = "Jamille Coelho Coimbra"
$ws.Range("R16").Value = 5111

# --- Row 17 ---------------------------------------------------------------
$ws.Range("A17").Value = 1261002320
$ws.Range("B17").Value = 3
$ws.Range("D17").Value = "BM"
$ws.Range("F17").Value = "NO"
$ws.Range("Z17").Value = "Desktop"
$ws.Range("N17").Value = "engenharia química bloco 3"
$ws.Range("I17").Value = "MONITOR LG 24 GAMER ULTRAGEAR FULL HD IPS 180HZ 1MS SRGB 99, HDR10, DISPLAYPORT HDMI G-SYNC FREESYNC, 24GS60F-B.AWZM EAN 789329995078"
$ws.Range("AA17").Value = "Jamille Coelho Coimbra"
$ws.Range("R17").Value = 5112

# Entered last so it lands as the final new shared string, matching the
# order the source workbook's table was built in.
$ws.Range("I16").Value = "Intel Core i9-14900KF 24-Core 32-Threads - Placa mãe Z790P - 128GB Memória DDR5 5600Mhz (4x32) - 2x SSD 1TB Nvme Kingston 4x4 L:7300Mbs G:6000Mbs - Placa de Video RTX 3060 NVIDIA 12GB 192Bits Cuda: 3584 - Fonte ATX 3.0 Modular 1250W 80Plus Gold PFC Ativo Cooler Master - Gabinete Masterbox MB520 - W11P"

# --- Formatting -------------------------------------------------------------
# The pasted "org/loc/pes" columns carry a Verdana 11 font (style s=2).
$verdanaRange = $ws.Range("N16,R16,AA16,N17,R17,AA17")
$verdanaRange.Font.Name = "Verdana"
$verdanaRange.Font.Size = 11
$verdanaRange.Font.Color = 0

# The "bem_cod/bem_dgv/bem_dsc_com" columns carry an Arial 10 font (style s=3).
$arialRange = $ws.Range("A16,B16,I16,A17,B17,I17")
$arialRange.Font.Name = "Arial"
$arialRange.Font.Size = 10

# Row heights for the two newly pasted rows.
$ws.Rows.Item(16).RowHeight = 14.25
$ws.Rows.Item(17).RowHeight = 14.25

# --- Selection / view --------------------------------------------------------
$ws.Range("A16:B17").Select()

$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
